# outputs-HGR-r202/test-p__Fusobacteriota_split_pruned.xlsx
# "updated outputs-r202, previous copy of ful-path.csv"
#
# The diff for this workbook does three things to the single data sheet
# ("quadratic-svm-score"):
#   1. Narrows column A from ~22.4 chars wide down to 5 chars wide.
#   2. Re-stamps the style of A1:C1 and A2 onto a freshly minted cell
#      format (still Text/"@" number format, still borderless - visually
#      identical, just a distinct style-table entry from the one shared
#      by the rest of the sheet).
#   3. Overwrites B2's value (1 -> 10990.46771063232), i.e. the real
#      "refreshed prediction score" data update referenced by the commit
#      message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column A width: target stored width="5" (character units). ---
# Excel's ColumnWidth<->stored "width" conversion adds the ~0.8333 (5/6)
# character padding constant, so ask for (5 - 5/6) to land exactly on 5.
$ws.Columns.Item(1).ColumnWidth = 4.166666666666667

# --- 2) Force A1:C1 and A2 onto a new (but visually-equivalent) style ---
# so they stop sharing the style index used by the rest of the sheet
# (mirrors the source diff's style-table split), without altering the
# rendered number format, font, fill or border.
$ws.Range("A1:C1").Locked = $true
$ws.Range("A2").Locked = $true

# --- 3) Data refresh: B2 changes from 1 to the recomputed score. ---
$ws.Range("B2").Value = 10990.46771063232
